$wb = $excel.ActiveWorkbook

$wsLanc = $wb.Worksheets.Item("LANÇAMENTOS")
$wsLanc.Range("A2:F4").ClearContents()

$wsEntrada = $wb.Worksheets.Item("ENTRADA")
$wsEntrada.Range("A2:D3").ClearContents()
